# Generate Report for Handoff
# Mark "b.md" as newly handed off (Ready for handoff) across the Overview,
# zh-cn and de-de sheets, with fresh handoff file names / timestamps.

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"

# --- Overview sheet: row for b.md is row 3 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status
$wsOverview.Range("D3").Value = "2016-03-23 08:37:43"

# --- zh-cn sheet: row for b.md is row 3 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $status
$wsZhCn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-23 08:37:38"

# --- de-de sheet: row for b.md is row 3 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $status
$wsDeDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-23 08:37:43"
